$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D (Price) updates — values that could be misread as a number
# get a leading apostrophe so Excel stores them as literal text,
# matching the inline-string cells already in the sheet.
$ws.Cells.Item(2, 4).Value = '69.994.15'
$ws.Cells.Item(3, 4).Value = '3.593.09'
$ws.Cells.Item(5, 4).Value = '''588.30'
$ws.Cells.Item(6, 4).Value = '''190.89'
$ws.Cells.Item(8, 4).Value = '3.583.94'
$ws.Cells.Item(10, 4).Value = '''0.177'
$ws.Cells.Item(11, 4).Value = '''0.659'
$ws.Cells.Item(12, 4).Value = '''57.88'
$ws.Cells.Item(15, 4).Value = '4.172.70'
$ws.Cells.Item(16, 4).Value = '3.601.91'
$ws.Cells.Item(18, 4).Value = '69.993.32'
$ws.Cells.Item(19, 4).Value = '''12.48'
$ws.Cells.Item(22, 4).Value = '''495.75'
$ws.Cells.Item(23, 4).Value = '''17.45'
$ws.Cells.Item(24, 4).Value = '''5.36'
$ws.Cells.Item(25, 4).Value = '''4.46'
$ws.Cells.Item(26, 4).Value = '''90.45'
$ws.Cells.Item(27, 4).Value = '''3.09'
$ws.Cells.Item(28, 4).Value = '''11.10'
$ws.Cells.Item(29, 4).Value = '''9.37'
$ws.Cells.Item(30, 4).Value = '''32.14'
$ws.Cells.Item(32, 4).Value = '''620.75'
$ws.Cells.Item(35, 4).Value = '''65.19'
$ws.Cells.Item(39, 4).Value = '''37.94'
$ws.Cells.Item(41, 4).Value = '''3.64'
$ws.Cells.Item(42, 4).Value = '3.314.81'
$ws.Cells.Item(44, 4).Value = '''0.0445'
$ws.Cells.Item(45, 4).Value = '''2.66'
$ws.Cells.Item(46, 4).Value = '''3.33'
$ws.Cells.Item(48, 4).Value = '''9.09'
$ws.Cells.Item(50, 4).Value = '''3.31'
$ws.Cells.Item(51, 4).Value = '''0.999'

# Column E (Volume 1h) updates
$ws.Cells.Item(2, 5).Value = '  +5.66%  '
$ws.Cells.Item(3, 5).Value = '  +5.28%  '
$ws.Cells.Item(4, 5).Value = '  -0.07%  '
$ws.Cells.Item(5, 5).Value = '  +3.54%  '
$ws.Cells.Item(6, 5).Value = '  +5.40%  '
$ws.Cells.Item(7, 5).Value = '  +2.14%  '
$ws.Cells.Item(8, 5).Value = '  +5.31%  '
$ws.Cells.Item(9, 5).Value = '  +0.07%  '
$ws.Cells.Item(10, 5).Value = '  -1.12%  '
$ws.Cells.Item(11, 5).Value = '  +2.67%  '
$ws.Cells.Item(12, 5).Value = '  +5.37%  '
$ws.Cells.Item(13, 5).Value = '  +4.01%  '
$ws.Cells.Item(14, 5).Value = '  +4.09%  '
$ws.Cells.Item(15, 5).Value = '  +4.97%  '
$ws.Cells.Item(16, 5).Value = '  +5.30%  '
$ws.Cells.Item(17, 5).Value = '  +5.34%  '
$ws.Cells.Item(18, 5).Value = '  +5.55%  '
$ws.Cells.Item(19, 5).Value = '  +3.99%  '
$ws.Cells.Item(20, 5).Value = '  +0.38%  '
$ws.Cells.Item(21, 5).Value = '  +3.99%  '
$ws.Cells.Item(22, 5).Value = '  +6.39%  '
$ws.Cells.Item(23, 5).Value = '  +19.60%  '
$ws.Cells.Item(24, 5).Value = '  +7.50%  '
$ws.Cells.Item(25, 5).Value = '  +7.34%  '
$ws.Cells.Item(26, 5).Value = '  +0.79%  '
$ws.Cells.Item(27, 5).Value = '  +5.40%  '
$ws.Cells.Item(28, 5).Value = '  +2.09%  '
$ws.Cells.Item(29, 5).Value = '  +5.69%  '
$ws.Cells.Item(30, 5).Value = '  +2.47%  '
$ws.Cells.Item(31, 5).Value = '  +8.34%  '
$ws.Cells.Item(32, 5).Value = '  +5.97%  '
$ws.Cells.Item(33, 5).Value = '  +5.21%  '
$ws.Cells.Item(34, 5).Value = '  +7.03%  '
$ws.Cells.Item(35, 5).Value = '  +4.31%  '
$ws.Cells.Item(36, 5).Value = '  +7.32%  '
$ws.Cells.Item(37, 5).Value = '  +5.41%  '
$ws.Cells.Item(38, 5).Value = '  +0.08%  '
$ws.Cells.Item(39, 5).Value = '  +4.33%  '
$ws.Cells.Item(40, 5).Value = '  +0.64%  '
$ws.Cells.Item(41, 5).Value = '  +1.10%  '
$ws.Cells.Item(42, 5).Value = '  +5.56%  '
$ws.Cells.Item(43, 5).Value = '  +4.52%  '
$ws.Cells.Item(44, 5).Value = '  +4.51%  '
$ws.Cells.Item(45, 5).Value = '  +5.45%  '
$ws.Cells.Item(46, 5).Value = '  +4.74%  '
$ws.Cells.Item(47, 5).Value = '  +1.84%  '
$ws.Cells.Item(48, 5).Value = '  +5.74%  '
$ws.Cells.Item(49, 5).Value = '  -3.76%  '
$ws.Cells.Item(50, 5).Value = '  +5.20%  '
$ws.Cells.Item(51, 5).Value = '  -0.20%  '

$wb.Save()
